$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert a new row at position 19 (shifts rows 19+ down by one, carries mergeCells along)
$ws.Rows("19:19").Insert()

# 2) Copy the ("last row") formatting that row 18 still has onto the brand-new row 19,
#    BEFORE we touch row 18's own formatting.
$ws.Range("B18:J18").Copy() | Out-Null
$ws.Range("B19:J19").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

# 3) Now re-style row 18 as an "interior" row (matching rows 16/17) since it is no
#    longer the last row of the table.
$ws.Range("B17:J17").Copy() | Out-Null
$ws.Range("B18:J18").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$excel.CutCopyMode = 0

# 4) Update VALOR MORA total and Cant. Periodos count
$ws.Range("E11").Value = 176940
$ws.Range("F13").Value = 4

# 5) Swap the "Novedad de Ingreso" / "Novedad de Retiro" header columns
$ws.Range("H15").Value = "Novedad de Retiro"
$ws.Range("I15").Value = "Novedad de Ingreso"

# 6) Reorder the period values (rows 16 & 18 swap; row 17 is unchanged) and refresh
#    the "Salario Basico" amounts for all existing data rows.
$ws.Range("E16").Value = "2204"
$ws.Range("G16").Value = 1423500

$ws.Range("G17").Value = 1423500

$ws.Range("E18").Value = "2206"
$ws.Range("G18").Value = 1423500

# 7) Populate the newly inserted row 19 with the additional 2508 period data
$ws.Range("B19").Value = "CC"
$ws.Range("C19").Value = "1007211195"
$ws.Range("D19").Value = "YORK DAVID CASTRO FUENTES"
$ws.Range("E19").Value = "2508"
$ws.Range("F19").Value = 56940
$ws.Range("G19").Value = 1423500
